# "Add test to Bull_Bear_beta macro"
# Append a new test-case row to the tests sheet: Test name / Description / macro name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").Value = "Bull/Bear beta1"
$ws.Range("B42").Value = "Test bull/bear beta"
$ws.Range("C42").Value = "Bull_Bear_beta_test1"

# Mirror the author's final view state (scrolled down, cursor parked on B45).
$ws.Range("B45").Select()
